# Update election result totals for BRAGA / GUIMARÃES (row 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1644
$ws.Range("I2").Value = 4352
$ws.Range("J2").Value = 17849
$ws.Range("K2").Value = 85
$ws.Range("L2").Value = 4919
$ws.Range("M2").Value = 277
$ws.Range("N2").Value = 3002
$ws.Range("O2").Value = 11
$ws.Range("P2").Value = 79
$ws.Range("Q2").Value = 24
$ws.Range("R2").Value = 239
$ws.Range("S2").Value = 1901
$ws.Range("T2").Value = 3164
$ws.Range("U2").Value = 236
$ws.Range("V2").Value = 27270
$ws.Range("W2").Value = 10
$ws.Range("X2").Value = 27399
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 408
$ws.Range("AA2").Value = 198
